$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "BMM"
$ws.Range("A7").Value = "Nissan"

$ws.Range("A8").Select()
